# Update the date title and the division problems in the practice table.
# The division sign "÷" (U+00F7) is built via [char] to avoid any
# script-encoding pitfalls.
$d = $word.ActiveDocument
$divide = [char]0x00F7

# --- Title line -------------------------------------------------------
$d.Content.Find.Execute("2024-01-30 Tuesday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-01-31 Wednesday", 2)

# --- Practice table -----------------------------------------------------
# The table is one 20-row x 5-column table; only every 4th row (1, 5, 9,
# 13, 17) holds problems, the rest are spacer rows.
$table = $d.Tables(1)

$grid = @(
    @{ Row = 1;  Cells = @("37${divide}7=", "75${divide}2=", "10${divide}3=", "87${divide}9=", "92${divide}7=") ; New = @("55${divide}8=", "43${divide}7=", "95${divide}8=", "87${divide}8=", "76${divide}7=") },
    @{ Row = 5;  Cells = @("44${divide}2=", "68${divide}3=", "19${divide}5=", "39${divide}5=", "33${divide}9=") ; New = @("33${divide}8=", "61${divide}8=", "81${divide}5=", "50${divide}6=", "77${divide}5=") },
    @{ Row = 9;  Cells = @("16${divide}7=", "75${divide}2=", "18${divide}3=", "43${divide}5=", "29${divide}2=") ; New = @("38${divide}4=", "57${divide}2=", "61${divide}5=", "60${divide}9=", "14${divide}3=") },
    @{ Row = 13; Cells = @("44${divide}9=", "57${divide}7=", "62${divide}4=", "68${divide}5=", "64${divide}3=") ; New = @("57${divide}2=", "80${divide}6=", "60${divide}3=", "60${divide}4=", "34${divide}5=") },
    @{ Row = 17; Cells = @("35${divide}4=", "67${divide}5=", "84${divide}6=", "38${divide}3=", "55${divide}7=") ; New = @("48${divide}8=", "16${divide}5=", "68${divide}4=", "38${divide}4=", "39${divide}2=") }
)

foreach ($entry in $grid) {
    $rowIndex = $entry.Row
    for ($col = 1; $col -le 5; $col++) {
        $expectedOld = $entry.Cells[$col - 1]
        $newValue = $entry.New[$col - 1]
        $cell = $table.Cell($rowIndex, $col)

        # Cell.Range.Text includes the trailing end-of-cell mark (CR+BEL);
        # trim it off before comparing against the plain expected string.
        $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)

        if ($currentText -ne $expectedOld) {
            Write-Output "WARNING: row $rowIndex col $col expected '$expectedOld' but found '$currentText'"
        }

        $cell.Range.Text = $newValue
    }
}
